$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 668.7826
$ws.Range("I2").Value = 248.42857
$ws.Range("J2").Value = 852.6875
$ws.Range("K2").Value = 248.42857
$ws.Range("L2").Value = 852.6875
$ws.Range("M2").Value = -135.42857
$ws.Range("N2").Value = -1078.6875
$ws.Range("H21").Value = 13261
$ws.Range("I21").Value = 7869.8
$ws.Range("K21").Value = 7869.8
$ws.Range("M21").Value = -7401.8
$ws.Range("H23").Value = 13261
$ws.Range("I23").Value = 7869.8
$ws.Range("K23").Value = 7869.8
$ws.Range("M23").Value = -7635.8
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1500
$ws.Range("K80").Value = 4500
$ws.Range("M80").Value = -3502
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1500
$ws.Range("K83").Value = 13500
$ws.Range("M83").Value = -8508
$ws.Range("H141").Value = 4242.5713
$ws.Range("I141").Value = 2449.6667
$ws.Range("K141").Value = 7349.000100000001
$ws.Range("M141").Value = -2169.000100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H61").Value = 7624.25
$ws.Range("I61").Value = 7624.25
$ws.Range("K61").Value = 7624.25
$ws.Range("M61").Value = -7412.25
$ws.Range("H74").Value = 3155
$ws.Range("I74").Value = 2835.5881
$ws.Range("J74").Value = 4965
$ws.Range("K74").Value = 2835.5881
$ws.Range("L74").Value = 4965
$ws.Range("M74").Value = -1961.5881
$ws.Range("N74").Value = -6713
$ws.Range("H77").Value = 3155
$ws.Range("I77").Value = 2835.5881
$ws.Range("J77").Value = 4965
$ws.Range("K77").Value = 14177.9405
$ws.Range("L77").Value = 24825
$ws.Range("M77").Value = -9809.940500000001
$ws.Range("N77").Value = -33561
$ws.Range("H110").Value = 4117263.5
$ws.Range("I110").Value = 5556935.5
$ws.Range("J110").Value = 3914.2856
$ws.Range("K110").Value = 5556935.5
$ws.Range("L110").Value = 3914.2856
$ws.Range("M110").Value = -5554890.5
$ws.Range("N110").Value = -8004.2856
$ws.Range("H122").Value = 748913.9
$ws.Range("I122").Value = 864649.5
$ws.Range("K122").Value = 2593948.5
$ws.Range("M122").Value = -2591498.5
$ws.Range("H136").Value = 7624.25
$ws.Range("I136").Value = 7624.25
$ws.Range("K136").Value = 22872.75
$ws.Range("M136").Value = -20322.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 100003380
$ws.Range("I16").Value = 100003380
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 100003380
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -100003093
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 3265
$ws.Range("I31").Value = 2083.3333
$ws.Range("K31").Value = 2083.3333
$ws.Range("M31").Value = -1788.3333
$ws.Range("H34").Value = 3265
$ws.Range("I34").Value = 2083.3333
$ws.Range("K34").Value = 2083.3333
$ws.Range("M34").Value = -1881.3333
$ws.Range("H99").Value = 11889.1
$ws.Range("I99").Value = 8425.923000000001
$ws.Range("K99").Value = 8425.923000000001
$ws.Range("M99").Value = -6927.923000000001
$ws.Range("H113").Value = 100003380
$ws.Range("I113").Value = 100003380
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 100003380
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -100001210
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 11889.1
$ws.Range("I126").Value = 8425.923000000001
$ws.Range("K126").Value = 25277.769
$ws.Range("M126").Value = -22807.769

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 43.125
$ws.Range("I10").Value = 35
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 105
$ws.Range("L10").Value = 300
$ws.Range("M10").Value = 34
$ws.Range("N10").Value = -578
$ws.Range("H34").Value = 4340
$ws.Range("J34").Value = 8333.333000000001
$ws.Range("L34").Value = 24999.999
$ws.Range("N34").Value = -25167.999
$ws.Range("H40").Value = 212.33333
$ws.Range("I40").Value = 212.33333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 849.33332
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -780.33332
$ws.Range("N40").ClearContents()
$ws.Range("H132").Value = 4284.5713
$ws.Range("J132").Value = 4249.75
$ws.Range("L132").Value = 38247.75
$ws.Range("N132").Value = -43307.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 49499.25
$ws.Range("J15").Value = 49499.25
$ws.Range("L15").Value = 49499.25
$ws.Range("N15").Value = -50075.25
$ws.Range("H80").Value = 3544.6
$ws.Range("I80").Value = 2851.7273
$ws.Range("J80").Value = 5450
$ws.Range("K80").Value = 2851.7273
$ws.Range("L80").Value = 5450
$ws.Range("M80").Value = -1853.7273
$ws.Range("N80").Value = -7446
$ws.Range("H81").Value = 49499.25
$ws.Range("J81").Value = 49499.25
$ws.Range("L81").Value = 49499.25
$ws.Range("N81").Value = -51495.25
$ws.Range("H83").Value = 3544.6
$ws.Range("I83").Value = 2851.7273
$ws.Range("J83").Value = 5450
$ws.Range("K83").Value = 14258.6365
$ws.Range("L83").Value = 27250
$ws.Range("M83").Value = -9266.636500000001
$ws.Range("N83").Value = -37234
$ws.Range("H84").Value = 49499.25
$ws.Range("J84").Value = 49499.25
$ws.Range("L84").Value = 148497.75
$ws.Range("N84").Value = -158481.75
$ws.Range("H107").Value = 1988.6666
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 4166
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 4166
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -8006
$ws.Range("H113").Value = 50028000
$ws.Range("I113").Value = 83358340
$ws.Range("K113").Value = 83358340
$ws.Range("M113").Value = -83356170

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15876958
$ws.Range("I61").Value = 22225822
$ws.Range("K61").Value = 22225822
$ws.Range("M61").Value = -22225620
$ws.Range("H82").Value = 47374.773
$ws.Range("I82").Value = 2324.75
$ws.Range("J82").Value = 167508.17
$ws.Range("K82").Value = 2324.75
$ws.Range("L82").Value = 167508.17
$ws.Range("M82").Value = -1963.75
$ws.Range("N82").Value = -168230.17
$ws.Range("H85").Value = 47374.773
$ws.Range("I85").Value = 2324.75
$ws.Range("J85").Value = 167508.17
$ws.Range("K85").Value = 2324.75
$ws.Range("L85").Value = 167508.17
$ws.Range("M85").Value = -1076.75
$ws.Range("N85").Value = -170004.17
$ws.Range("H113").Value = 15876958
$ws.Range("I113").Value = 22225822
$ws.Range("K113").Value = 22225822
$ws.Range("M113").Value = -22223652

Write-Host "Applied 171 sets and 4 clears"
